$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Sheet1
$ws2 = $wb.Worksheets.Item(2)   # foo
$ws3 = $wb.Worksheets.Item(3)   # Sheet3 -> bar

# --- Rename Sheet3 -> bar -------------------------------------------------
$ws3.Name = "bar"

# --- Rebuild the pre-existing external link ("external.xlsx") from scratch
#     so its sheet-name table becomes Sheet1, Sheet2, Sheet3, bar (no stale
#     blank entry) while keeping the formula's [1] link-index notation.
$wb.BreakLink("external.xlsx", 1)
$ws2.Range("A1").ClearContents()

$ws2.Range("Z1").Formula = "=[external.xlsx]Sheet1!A1"
$ws2.Range("Z2").Formula = "=[external.xlsx]Sheet2!A1"
$ws2.Range("Z3").Formula = "=[external.xlsx]Sheet3!A1"
$ws2.Range("D7").Formula = "=[external.xlsx]bar!C7"
$ws2.Range("D7").Formula = "=+[1]bar!`$C`$7"
$ws2.Range("Z1:Z3").ClearContents()

# --- Second external link ("external2.xlsx"): Sheet1, bar ----------------
$ws2.Range("Z1").Formula = "=[external2.xlsx]Sheet1!A1"
$ws2.Range("D8").Formula = "=[external2.xlsx]bar!C8"
$ws2.Range("D8").Formula = "=+[2]bar!`$C`$8"
$ws2.Range("Z1").ClearContents()

# --- Third external link ("external3.xlsx"): bar --------------------------
$ws2.Range("D9").Formula = "=[external3.xlsx]bar!C9"
$ws2.Range("D9").Formula = "=+[3]bar!`$C`$9"

# D7 carries a percentage number format (new cellXfs entry numFmtId=9).
$ws2.Range("D7").NumberFormat = "0%"

# --- sheet "bar" (was Sheet3) content --------------------------------------
$ws3.Range("A1").Value = 1
$ws3.Range("A2").Value = 2
$ws3.Range("A3").Value = 3
$ws3.Range("C7").Formula = "=SUM(A1:A3)"
$ws3.Range("C9").Value = "text"

# --- selections / active sheet --------------------------------------------
$ws2.Activate()
$ws2.Range("D9").Select()
$ws3.Activate()
$ws3.Range("C9").Select()
